# Applies the Dec 28 2023 cryptos-list refresh: updated prices/
# volume percentages across most rows, plus a name/link/price swap
# between the InjectiveProtocol and Monero rows (29 <-> 30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.004.23"
$ws.Range("E2").Value = "  +0.99%  "
# Row 3
$ws.Range("D3").Value = "2.372.20"
$ws.Range("E3").Value = "  +6.56%  "
# Row 4
$ws.Range("E4").Value = "  -0.36%  "
# Row 5
$ws.Range("D5").Value = "'322.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.30%  "
# Row 6
$ws.Range("D6").Value = "'103.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.60%  "
# Row 7
$ws.Range("E7").Value = "  +2.25%  "
# Row 8
$ws.Range("E8").Value = "  +0.02%  "
# Row 9
$ws.Range("D9").Value = "'0.655"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.69%  "
# Row 10
$ws.Range("D10").Value = "'41.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.88%  "
# Row 11
$ws.Range("D11").Value = "'0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.62%  "
# Row 12
$ws.Range("D12").Value = "'8.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.36%  "
# Row 13
$ws.Range("E13").Value = "  -2.63%  "
# Row 14
$ws.Range("D14").Value = "'16.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.43%  "
# Row 15
$ws.Range("E15").Value = "  +2.07%  "
# Row 16
$ws.Range("D16").Value = "2.736.68"
$ws.Range("E16").Value = "  +6.87%  "
# Row 17
$ws.Range("D17").Value = "2.375.49"
$ws.Range("E17").Value = "  +6.88%  "
# Row 18
$ws.Range("D18").Value = "43.076.78"
$ws.Range("E18").Value = "  +1.53%  "
# Row 19
$ws.Range("D19").Value = "'7.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.57%  "
# Row 20
$ws.Range("E20").Value = "  +2.09%  "
# Row 21
$ws.Range("D21").Value = "'76.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.07%  "
# Row 22
$ws.Range("D22").Value = "'275.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.32%  "
# Row 23
$ws.Range("D23").Value = "'3.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
# Row 24
$ws.Range("E24").Value = "  +1.44%  "
# Row 25
$ws.Range("D25").Value = "'9.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.68%  "
# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "
# Row 27
$ws.Range("E27").Value = "  +1.89%  "
# Row 28
$ws.Range("D28").Value = "'23.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.56%  "
# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'175.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'37.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "
# Row 31
$ws.Range("E31").Value = "  -1.44%  "
# Row 32
$ws.Range("E32").Value = "  +1.98%  "
# Row 33
$ws.Range("E33").Value = "  +4.57%  "
# Row 34
$ws.Range("E34").Value = "  +2.66%  "
# Row 35
$ws.Range("D35").Value = "'0.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.92%  "
# Row 36
$ws.Range("D36").Value = "'4.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.21%  "
# Row 37
$ws.Range("D37").Value = "'4.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "
# Row 38
$ws.Range("E38").Value = "  -2.15%  "
# Row 39
$ws.Range("E39").Value = "  +1.46%  "
# Row 40
$ws.Range("D40").Value = "'2.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.63%  "
# Row 41
$ws.Range("E41").Value = "  +21.69%  "
# Row 42
$ws.Range("D42").Value = "'123.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +21.49%  "
# Row 43
$ws.Range("D43").Value = "'0.230"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.91%  "
# Row 44
$ws.Range("D44").Value = "'69.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.95%  "
# Row 45
$ws.Range("E45").Value = "  +0.28%  "
# Row 46
$ws.Range("D46").Value = "'93.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +60.51%  "
# Row 47
$ws.Range("D47").Value = "'12.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.74%  "
# Row 48
$ws.Range("D48").Value = "'9.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.53%  "
# Row 49
$ws.Range("D49").Value = "'5.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.74%  "
# Row 50
$ws.Range("E50").Value = "  +1.02%  "
# Row 51
$ws.Range("D51").Value = "1.600.87"
$ws.Range("E51").Value = "  +12.05%  "
